$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.147999999999989
$ws.Range("A9").Value = -20.52319999999998
$ws.Range("A18").Value = -23.04250000000002
$ws.Range("A20").Value = -22.10340000000003
$ws.Range("C21").Value = -13.2265
